# Variables sheet: rename the METHOD column to DERIVATION, and insert a new
# empty METHOD column right after it (kept blank, for backward reference).
#
# Column layout before: ... Q=VLM_ID, R=METHOD,     S=DEPENDS_ON
# Column layout after:  ... Q=VLM_ID, R=DERIVATION, S=METHOD (blank), T=DEPENDS_ON

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Variables")

# Insert a new blank column at S. This pushes the existing DEPENDS_ON
# column (S) out to T, while column R (METHOD data) is left untouched.
$ws.Range("S1").EntireColumn.Insert()

# Column R's data stays as-is; just re-label its header from METHOD to
# DERIVATION (renaming the column in place).
$ws.Range("R1").Value = "DERIVATION"

# The freshly inserted, empty column S becomes the new (blank) METHOD
# column, kept for backward reference.
$ws.Range("S1").Value = "METHOD"
